# Weekly update: a new price-report week was inserted for "Fruta, Terminal
# Hortofrutícola Agro Chillán - Palta". Two new rows of data (Hass, Primera /
# Segunda, fecha 44516) are inserted right before the old row 325, pushing
# every subsequent row down by two (the last two rows of the sheet become
# brand-new rows as a consequence of that shift).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 325:326 - this shifts rows 325.. down to 327..
$ws.Rows("325:326").Insert()

# --- Row 325 (new) ---------------------------------------------------
$ws.Cells.Item(325,1).Value  = 7
$ws.Cells.Item(325,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(325,3).Value  = "Ñuble"
$ws.Cells.Item(325,4).Value  = 44516
$ws.Cells.Item(325,5).Value  = 16
$ws.Cells.Item(325,6).Value  = "Fruta"
$ws.Cells.Item(325,7).Value  = 100106
$ws.Cells.Item(325,8).Value  = "Oleaginosos"
$ws.Cells.Item(325,9).Value  = 100106002
$ws.Cells.Item(325,10).Value = "Palta"
$ws.Cells.Item(325,11).Value = "Hass"
$ws.Cells.Item(325,12).Value = "Primera"
$ws.Cells.Item(325,13).Value = 60
$ws.Cells.Item(325,14).Value = 2700
$ws.Cells.Item(325,15).Value = 2800
$ws.Cells.Item(325,16).Value = 2750
$ws.Cells.Item(325,17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(325,18).Value = "Provincia de Quillota"
$ws.Cells.Item(325,19).Value = 2750
$ws.Cells.Item(325,20).Value = 1

# --- Row 326 (new) ---------------------------------------------------
$ws.Cells.Item(326,1).Value  = 7
$ws.Cells.Item(326,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(326,3).Value  = "Ñuble"
$ws.Cells.Item(326,4).Value  = 44516
$ws.Cells.Item(326,5).Value  = 16
$ws.Cells.Item(326,6).Value  = "Fruta"
$ws.Cells.Item(326,7).Value  = 100106
$ws.Cells.Item(326,8).Value  = "Oleaginosos"
$ws.Cells.Item(326,9).Value  = 100106002
$ws.Cells.Item(326,10).Value = "Palta"
$ws.Cells.Item(326,11).Value = "Hass"
$ws.Cells.Item(326,12).Value = "Segunda"
$ws.Cells.Item(326,13).Value = 120
$ws.Cells.Item(326,14).Value = 2400
$ws.Cells.Item(326,15).Value = 2500
$ws.Cells.Item(326,16).Value = 2450
$ws.Cells.Item(326,17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(326,18).Value = "Provincia de Quillota"
$ws.Cells.Item(326,19).Value = 2450
$ws.Cells.Item(326,20).Value = 1

# Keep the date column's number format consistent with the rest of the
# column (Insert() already carried style "2" down onto D325/D326, but set
# it explicitly too so it is robust regardless of engine defaults).
$ws.Range("D325:D326").NumberFormat = $ws.Range("D324").NumberFormat
